$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.144.86'
$ws.Range('E2').Value = '  -1.83%  '
$ws.Range('D3').Value = '2.245.81'
$ws.Range('E3').Value = '  -1.83%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = '''246.67'
$ws.Range('E5').Value = '  -2.20%  '
$ws.Range('E6').Value = '  +0.45%  '
$ws.Range('D7').Value = '''76.69'
$ws.Range('E7').Value = '  +4.44%  '
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('D9').Value = '''0.624'
$ws.Range('E9').Value = '  -4.09%  '
$ws.Range('D10').Value = '''41.19'
$ws.Range('E10').Value = '  +5.33%  '
$ws.Range('D11').Value = '''0.0954'
$ws.Range('E11').Value = '  -2.29%  '
$ws.Range('E12').Value = '  -4.47%  '
$ws.Range('D13').Value = '''0.102'
$ws.Range('E13').Value = '  -3.19%  '
$ws.Range('D14').Value = '2.581.01'
$ws.Range('E14').Value = '  -1.85%  '
$ws.Range('D15').Value = '''14.78'
$ws.Range('E15').Value = '  -2.96%  '
$ws.Range('D16').Value = '''0.860'
$ws.Range('E16').Value = '  -1.07%  '
$ws.Range('D17').Value = '2.238.10'
$ws.Range('E17').Value = '  -2.11%  '
$ws.Range('D18').Value = '42.034.05'
$ws.Range('E18').Value = '  -1.87%  '
$ws.Range('D19').Value = '0.0₃0984'
$ws.Range('E19').Value = '  -2.32%  '
$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D20').Value = '''6.11'
$ws.Range('E20').Value = '  -2.77%  '
$ws.Range('B21').Value = 'Litecoin'
$ws.Range('C21').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D21').Value = '''71.87'
$ws.Range('E21').Value = '  -0.90%  '
$ws.Range('D22').Value = '''2.30'
$ws.Range('E22').Value = '  +2.54%  '
$ws.Range('D23').Value = '''230.26'
$ws.Range('E23').Value = '  -2.76%  '
$ws.Range('E24').Value = '  +0.00%  '
$ws.Range('D25').Value = '''11.25'
$ws.Range('E25').Value = '  -2.66%  '
$ws.Range('D26').Value = '''3.62'
$ws.Range('E26').Value = '  -7.37%  '
$ws.Range('D27').Value = '''2.29'
$ws.Range('E27').Value = '  -4.77%  '
$ws.Range('D28').Value = '''7.36'
$ws.Range('E28').Value = '  +14.05%  '
$ws.Range('E29').Value = '  +0.94%  '
$ws.Range('D30').Value = '''169.46'
$ws.Range('E30').Value = '  +1.50%  '
$ws.Range('D31').Value = '''20.60'
$ws.Range('E31').Value = '  -2.11%  '
$ws.Range('D32').Value = '''33.11'
$ws.Range('E32').Value = '  +7.24%  '
$ws.Range('E33').Value = '  +0.09%  '
$ws.Range('D34').Value = '''0.120'
$ws.Range('E34').Value = '  -5.02%  '
$ws.Range('E35').Value = '  -1.13%  '
$ws.Range('E36').Value = '  -1.27%  '
$ws.Range('E37').Value = '  +2.85%  '
$ws.Range('D38').Value = '''0.0303'
$ws.Range('E38').Value = '  -1.70%  '
$ws.Range('D39').Value = '''14.07'
$ws.Range('E39').Value = '  +0.28%  '
$ws.Range('E40').Value = '  -0.16%  '
$ws.Range('E41').Value = '  -6.64%  '
$ws.Range('E42').Value = '  +11.40%  '
$ws.Range('E43').Value = '  -6.33%  '
$ws.Range('D44').Value = '''61.08'
$ws.Range('E44').Value = '  -1.07%  '
$ws.Range('E45').Value = '  -5.48%  '
$ws.Range('D46').Value = '''0.0999'
$ws.Range('E46').Value = '  -3.83%  '
$ws.Range('D47').Value = '''0.997'
$ws.Range('E47').Value = '  -0.38%  '
$ws.Range('E48').Value = '  -2.97%  '
$ws.Range('E49').Value = '  -1.00%  '
$ws.Range('E50').Value = '  +14.52%  '
$ws.Range('D51').Value = '''2.30'
$ws.Range('E51').Value = '  +0.55%  '
